# Weekly driver report update for 2025-04-28
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers summary (row 3 / totals row 4) ---
$ws.Range("C3").Value = 76
$ws.Range("D3").Value = 97.5
$ws.Range("C4").Value = 76

# --- Good Drivers table: a new driver version has appeared at the top
#     (21.40.1.3), pushing the rest of the list down by one row ---
$ws.Rows.Item(12).Insert()

# New row 12: freshly observed driver, no prior "Driver Vintage" date yet
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B12").Value = 11128
$ws.Range("B12").NumberFormat = "#,##0"
$ws.Range("D12").Value = 100

# Updated sample counts / dates for the rest of the (shifted) rows
$ws.Range("B13").Value = 486214
$ws.Range("B14").Value = 79953
$ws.Range("B15").Value = 35355
$ws.Range("B16").Value = 65425
$ws.Range("B17").Value = 117653

# Touch the bottom-right corner of the sheet's padding area so the
# worksheet's used range (and thus its dimension) keeps extending down
# to row 23 / column J, same as the rest of this weekly report template.
$ws.Range("J23").Borders.LineStyle = -4142
